$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1417.25
$ws.Range("I19").Value = 1417.25
$ws.Range("K19").Value = 1417.25
$ws.Range("M19").Value = -1242.25

$ws.Range("H53").Value = 370.7
$ws.Range("I53").Value = 370.7
$ws.Range("K53").Value = 370.7
$ws.Range("M53").Value = 266.3

$ws.Range("H112").Value = 1940.7
$ws.Range("I112").Value = 1245.8334
$ws.Range("J112").Value = 2238.5
$ws.Range("K112").Value = 3737.5002
$ws.Range("L112").Value = 6715.5
$ws.Range("M112").Value = -2629.5002
$ws.Range("N112").Value = -8931.5

$ws.Range("H135").Value = 1010.2857
$ws.Range("I135").Value = 1010.2857
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9092.5713
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6557.5713
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3786.7805
$ws.Range("I63").Value = 3284.3713
$ws.Range("K63").Value = 3284.3713
$ws.Range("M63").Value = -2598.3713

$ws.Range("H66").Value = 3786.7805
$ws.Range("I66").Value = 3284.3713
$ws.Range("K66").Value = 16421.8565
$ws.Range("M66").Value = -12989.8565

$ws.Range("H74").Value = 3320.0527
$ws.Range("J74").Value = 4458.143
$ws.Range("L74").Value = 4458.143
$ws.Range("N74").Value = -6206.143

$ws.Range("H77").Value = 3320.0527
$ws.Range("J77").Value = 4458.143
$ws.Range("L77").Value = 22290.715
$ws.Range("N77").Value = -31026.715

$ws.Range("H122").Value = 3402
$ws.Range("I122").Value = 2577.25
$ws.Range("K122").Value = 7731.75
$ws.Range("M122").Value = -5281.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 79759.91
$ws.Range("J132").Value = 79759.91
$ws.Range("L132").Value = 79759.91
$ws.Range("N132").Value = -89879.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3019.152
$ws.Range("I31").Value = 1302.875
$ws.Range("J31").Value = 3380.4736
$ws.Range("K31").Value = 1302.875
$ws.Range("L31").Value = 3380.4736
$ws.Range("M31").Value = -1007.875
$ws.Range("N31").Value = -3970.4736

$ws.Range("H34").Value = 3019.152
$ws.Range("I34").Value = 1302.875
$ws.Range("J34").Value = 3380.4736
$ws.Range("K34").Value = 1302.875
$ws.Range("L34").Value = 3380.4736
$ws.Range("M34").Value = -1100.875
$ws.Range("N34").Value = -3784.4736

$ws.Range("H58").Value = 5894.95
$ws.Range("I58").Value = 4159.294
$ws.Range("J58").Value = 15730.333
$ws.Range("K58").Value = 4159.294
$ws.Range("L58").Value = 15730.333
$ws.Range("M58").Value = -3956.294
$ws.Range("N58").Value = -16136.333

$ws.Range("H62").Value = 58117.79
$ws.Range("I62").Value = 94175.37
$ws.Range("K62").Value = 94175.37
$ws.Range("M62").Value = -93551.37

$ws.Range("H65").Value = 58117.79
$ws.Range("I65").Value = 94175.37
$ws.Range("K65").Value = 470876.85
$ws.Range("M65").Value = -467756.85

$ws.Range("H132").Value = 33441.098
$ws.Range("I132").Value = 24302.738
$ws.Range("K132").Value = 72908.21400000001
$ws.Range("M132").Value = -70378.21400000001

$ws.Range("H136").Value = 5894.95
$ws.Range("I136").Value = 4159.294
$ws.Range("J136").Value = 15730.333
$ws.Range("K136").Value = 12477.882
$ws.Range("L136").Value = 47190.999
$ws.Range("M136").Value = -9927.882
$ws.Range("N136").Value = -52290.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4989
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4989
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 14967
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -16465

$ws.Range("H66").Value = 4989
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4989
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 44901
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -52389

$ws.Range("H107").Value = 678
$ws.Range("I107").Value = 525.63635
$ws.Range("K107").Value = 1576.90905
$ws.Range("M107").Value = 343.09095

$ws.Range("H117").Value = 1109.75
$ws.Range("I117").Value = 977.4
$ws.Range("K117").Value = 2932.2
$ws.Range("M117").Value = 509.8000000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3287.25
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 3716.3333
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 11148.9999
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -16048.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 750
$ws.Range("I9").Value = 212.5
$ws.Range("K9").Value = 212.5
$ws.Range("M9").Value = 11.5

$ws.Range("H61").Value = 848.087
$ws.Range("I61").Value = 872.8095
$ws.Range("K61").Value = 872.8095
$ws.Range("M61").Value = -670.8095

$ws.Range("H113").Value = 848.087
$ws.Range("I113").Value = 872.8095
$ws.Range("K113").Value = 872.8095
$ws.Range("M113").Value = 1297.1905

$ws.Range("H122").Value = 2023.6
$ws.Range("I122").Value = 1492.2858
$ws.Range("J122").Value = 3263.3333
$ws.Range("K122").Value = 4476.857400000001
$ws.Range("L122").Value = 9789.999899999999
$ws.Range("M122").Value = -2026.857400000001
$ws.Range("N122").Value = -14689.9999

$ws.Range("H132").Value = 4332.84
$ws.Range("I132").Value = 3941.2
$ws.Range("J132").Value = 4920.3
$ws.Range("K132").Value = 11823.6
$ws.Range("L132").Value = 14760.9
$ws.Range("M132").Value = -9293.599999999999
$ws.Range("N132").Value = -19820.9

$ws.Range("H136").Value = 2726.1133
$ws.Range("I136").Value = 2429.422
$ws.Range("J136").Value = 4395
$ws.Range("K136").Value = 7288.266
$ws.Range("L136").Value = 13185
$ws.Range("M136").Value = -4738.266
$ws.Range("N136").Value = -18285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 36057.6
$ws.Range("J75").Value = 26129.666
$ws.Range("L75").Value = 26129.666
$ws.Range("N75").Value = -28001.666

$ws.Range("H78").Value = 36057.6
$ws.Range("J78").Value = 26129.666
$ws.Range("L78").Value = 78388.99800000001
$ws.Range("N78").Value = -87748.99800000001

$ws.Range("H122").Value = 6045.231
$ws.Range("I122").Value = 2658.9
$ws.Range("J122").Value = 17333
$ws.Range("K122").Value = 7976.700000000001
$ws.Range("L122").Value = 51999
$ws.Range("M122").Value = -5526.700000000001
$ws.Range("N122").Value = -56899
